$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '67.023.52'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -5.35%  '
$ws.Cells.Item(2, 5).ClearFormats()

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.366.42'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -6.81%  '
$ws.Cells.Item(3, 5).ClearFormats()

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(4, 5).ClearFormats()

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '562.04'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -6.35%  '
$ws.Cells.Item(5, 5).ClearFormats()

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '185.45'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -8.52%  '
$ws.Cells.Item(6, 5).ClearFormats()

# Row 7
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(7, 3).ClearFormats()
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.596'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  -5.19%  '
$ws.Cells.Item(7, 5).ClearFormats()

# Row 8
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(8, 3).ClearFormats()
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  +0.16%  '
$ws.Cells.Item(8, 5).ClearFormats()

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '3.357.19'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -6.73%  '
$ws.Cells.Item(9, 5).ClearFormats()

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.193'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -10.80%  '
$ws.Cells.Item(10, 5).ClearFormats()

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.598'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -7.48%  '
$ws.Cells.Item(11, 5).ClearFormats()

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '48.22'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -10.76%  '
$ws.Cells.Item(12, 5).ClearFormats()

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000270'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -11.05%  '
$ws.Cells.Item(13, 5).ClearFormats()

# Row 14
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 2).ClearFormats()
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 3).ClearFormats()
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '8.80'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  -8.60%  '
$ws.Cells.Item(14, 5).ClearFormats()

# Row 15
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 2).ClearFormats()
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 3).ClearFormats()
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.914.86'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -6.38%  '
$ws.Cells.Item(15, 5).ClearFormats()

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '614.30'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -10.18%  '
$ws.Cells.Item(16, 5).ClearFormats()

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.420.19'
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -5.01%  '
$ws.Cells.Item(17, 5).ClearFormats()

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '66.995.71'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -5.45%  '
$ws.Cells.Item(18, 5).ClearFormats()

# Row 19
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -3.29%  '
$ws.Cells.Item(19, 5).ClearFormats()

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.75'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  -7.70%  '
$ws.Cells.Item(20, 5).ClearFormats()

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.80'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -8.21%  '
$ws.Cells.Item(21, 5).ClearFormats()

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.921'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -7.97%  '
$ws.Cells.Item(22, 5).ClearFormats()

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '17.45'
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -7.44%  '
$ws.Cells.Item(23, 5).ClearFormats()

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '5.14'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -4.12%  '
$ws.Cells.Item(24, 5).ClearFormats()

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '95.43'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  -13.48%  '
$ws.Cells.Item(25, 5).ClearFormats()

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '4.13'
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -10.82%  '
$ws.Cells.Item(26, 5).ClearFormats()

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.77'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -9.17%  '
$ws.Cells.Item(27, 5).ClearFormats()

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.65'
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -9.21%  '
$ws.Cells.Item(28, 5).ClearFormats()

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.88'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -12.60%  '
$ws.Cells.Item(29, 5).ClearFormats()

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '31.02'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -10.00%  '
$ws.Cells.Item(30, 5).ClearFormats()

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.92'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -12.30%  '
$ws.Cells.Item(31, 5).ClearFormats()

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.38'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -11.06%  '
$ws.Cells.Item(32, 5).ClearFormats()

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '11.25'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  -8.60%  '
$ws.Cells.Item(33, 5).ClearFormats()

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '553.57'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  +7.30%  '
$ws.Cells.Item(34, 5).ClearFormats()

# Row 35
$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 2).ClearFormats()
$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(35, 3).ClearFormats()
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.106'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -8.12%  '
$ws.Cells.Item(35, 5).ClearFormats()

# Row 36
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = 'OKB'
$ws.Cells.Item(36, 2).ClearFormats()
$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(36, 3).ClearFormats()
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '59.10'
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -7.05%  '
$ws.Cells.Item(36, 5).ClearFormats()

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.809.12'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  -1.94%  '
$ws.Cells.Item(37, 5).ClearFormats()

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.996'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  -0.37%  '
$ws.Cells.Item(38, 5).ClearFormats()

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.76'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  +35.40%  '
$ws.Cells.Item(39, 5).ClearFormats()

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.47'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -3.69%  '
$ws.Cells.Item(40, 5).ClearFormats()

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0₃0736'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -14.20%  '
$ws.Cells.Item(41, 5).ClearFormats()

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.77'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -8.74%  '
$ws.Cells.Item(42, 5).ClearFormats()

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.357'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -7.20%  '
$ws.Cells.Item(43, 5).ClearFormats()

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.128'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -8.23%  '
$ws.Cells.Item(44, 5).ClearFormats()

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '32.72'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -11.59%  '
$ws.Cells.Item(45, 5).ClearFormats()

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0423'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  -9.82%  '
$ws.Cells.Item(46, 5).ClearFormats()

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.25'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -4.66%  '
$ws.Cells.Item(47, 5).ClearFormats()

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.72'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -11.30%  '
$ws.Cells.Item(48, 5).ClearFormats()

# Row 49
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -7.82%  '
$ws.Cells.Item(49, 5).ClearFormats()

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.996'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -0.56%  '
$ws.Cells.Item(50, 5).ClearFormats()

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.78'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -10.12%  '
$ws.Cells.Item(51, 5).ClearFormats()
